$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "69.243.77"
Set-TextValue "E2" "  +1.67%  "

# Row 3
Set-TextValue "D3" "3.788.01"
Set-TextValue "E3" "  -0.20%  "

# Row 4
Set-TextValue "E4" "  +0.17%  "

# Row 5
Set-TextValue "D5" "627.63"
Set-TextValue "E5" "  +4.37%  "

# Row 6
Set-TextValue "D6" "164.79"
Set-TextValue "E6" "  -0.23%  "

# Row 7
Set-TextValue "D7" "3.786.07"
Set-TextValue "E7" "  -0.19%  "

# Row 8
Set-TextValue "E8" "  -0.20%  "

# Row 9
Set-TextValue "D9" "0.521"
Set-TextValue "E9" "  +0.61%  "

# Row 10
Set-TextValue "E10" "  +1.27%  "

# Row 11
Set-TextValue "E11" "  +0.24%  "

# Row 12
Set-TextValue "D12" "6.67"
Set-TextValue "E12" "  +2.75%  "

# Row 13
Set-TextValue "E13" "  -0.43%  "

# Row 14
Set-TextValue "D14" "35.61"
Set-TextValue "E14" "  -0.64%  "

# Row 15
Set-TextValue "D15" "4.424.97"
Set-TextValue "E15" "  -0.15%  "

# Row 16
Set-TextValue "D16" "3.799.07"
Set-TextValue "E16" "  +0.88%  "

# Row 17
Set-TextValue "D17" "69.280.36"
Set-TextValue "E17" "  +1.70%  "

# Row 18
Set-TextValue "D18" "17.95"
Set-TextValue "E18" "  -2.25%  "

# Row 19
Set-TextValue "D19" "7.12"
Set-TextValue "E19" "  +0.62%  "

# Row 20
Set-TextValue "E20" "  -1.26%  "

# Row 21
Set-TextValue "D21" "468.80"
Set-TextValue "E21" "  +1.59%  "

# Row 22
Set-TextValue "D22" "9.64"
Set-TextValue "E22" "  -0.64%  "

# Row 23
Set-TextValue "E23" "  +0.40%  "

# Row 24
Set-TextValue "D24" "0.0000151"
Set-TextValue "E24" "  +1.93%  "

# Row 25
Set-TextValue "D25" "83.31"
Set-TextValue "E25" "  +0.33%  "

# Row 26
Set-TextValue "D26" "12.07"
Set-TextValue "E26" "  +0.43%  "

# Row 27
Set-TextValue "D27" "2.17"
Set-TextValue "E27" "  +2.38%  "

# Row 28
Set-TextValue "D28" "10.02"
Set-TextValue "E28" "  +0.09%  "

# Row 29
Set-TextValue "E29" "  -0.15%  "

# Row 30
Set-TextValue "D30" "3.935.22"
Set-TextValue "E30" "  -0.21%  "

# Row 31
Set-TextValue "E31" "  +0.98%  "

# Row 32
Set-TextValue "E32" "  +0.11%  "

# Row 33
Set-TextValue "E33" "  -0.37%  "

# Row 34
Set-TextValue "D34" "28.96"
Set-TextValue "E34" "  -1.28%  "

# Row 35
Set-TextValue "E35" "  -0.11%  "

# Row 36
Set-TextValue "E36" "  -0.10%  "

# Row 37
Set-TextValue "D37" "3.737.93"

# Row 38
Set-TextValue "E38" "  +2.70%  "

# Row 39
Set-TextValue "D39" "0.151"
Set-TextValue "E39" "  +8.55%  "

# Row 40
Set-TextValue "D40" "3.36"
Set-TextValue "E40" "  +0.48%  "

# Row 41
Set-TextValue "E41" "  -0.17%  "

# Row 42
Set-TextValue "D42" "0.967"
Set-TextValue "E42" "  -1.98%  "

# Row 43
Set-TextValue "D43" "1.00"
Set-TextValue "E43" "  +0.05%  "

# Row 45
Set-TextValue "E45" "  +0.20%  "

# Row 46
Set-TextValue "B46" "Stacks"
Set-TextValue "C46" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D46" "1.94"
Set-TextValue "E46" "  +3.58%  "

# Row 47
Set-TextValue "B47" "Monero"
Set-TextValue "C47" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D47" "153.19"
Set-TextValue "E47" "  +0.44%  "

# Row 48
Set-TextValue "D48" "46.84"
Set-TextValue "E48" "  -1.65%  "

# Row 49
Set-TextValue "D49" "42.85"
Set-TextValue "E49" "  -0.83%  "

# Row 50
Set-TextValue "B50" "ONDO"
Set-TextValue "C50" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D50" "1.39"
Set-TextValue "E50" "  +3.20%  "

# Row 51
Set-TextValue "B51" "Cosmos"
Set-TextValue "C51" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "8.44"
Set-TextValue "E51" "  +0.98%  "
